# Generate Report for handback
# The localization pipeline re-ran: the handback transform for
# 0479b554-9e75-4834-8f00-0baa74747d98.md failed, while
# 27b89366-4817-4392-9fec-e9a51b422d94.md is still in translation.
# Update the status report (swap the two rows' data and set the new status).

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A4").Value = "0479b554-9e75-4834-8f00-0baa74747d98.md"
$ws.Range("B4").Value = "Handback transform failed"
$ws.Range("C4").Value = "Handback transform failed"

$ws.Range("A5").Value = "27b89366-4817-4392-9fec-e9a51b422d94.md"
$ws.Range("B5").Value = "In Translation"
$ws.Range("C5").Value = "In Translation"

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("A4").Value = "0479b554-9e75-4834-8f00-0baa74747d98.md"
$ws.Range("B4").Value = "Handback transform failed"
$ws.Range("C4").Value = "0479b554-9e75-4834-8f00-0baa74747d98.42cc1743e91fdb6c13d01627e17a3bebdf93eac3.zh-cn.xlf"
$ws.Range("D4").Value = "2016-01-26 10:21:27"

$ws.Range("A5").Value = "27b89366-4817-4392-9fec-e9a51b422d94.md"
$ws.Range("B5").Value = "In Translation"
$ws.Range("C5").Value = "27b89366-4817-4392-9fec-e9a51b422d94.e41ac49300ef89bcc58f59aeef53d4d36f3c06ba.zh-cn.xlf"
$ws.Range("D5").Value = "2016-01-26 10:19:11"

# ---- de-de sheet ----
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("A4").Value = "0479b554-9e75-4834-8f00-0baa74747d98.md"
$ws.Range("B4").Value = "Handback transform failed"
$ws.Range("C4").Value = "0479b554-9e75-4834-8f00-0baa74747d98.42cc1743e91fdb6c13d01627e17a3bebdf93eac3.de-de.xlf"
$ws.Range("D4").Value = "2016-01-26 10:21:40"

$ws.Range("A5").Value = "27b89366-4817-4392-9fec-e9a51b422d94.md"
$ws.Range("B5").Value = "In Translation"
$ws.Range("C5").Value = "27b89366-4817-4392-9fec-e9a51b422d94.e41ac49300ef89bcc58f59aeef53d4d36f3c06ba.de-de.xlf"
$ws.Range("D5").Value = "2016-01-26 10:19:23"
